$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.171.74'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '1.682.66'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.58'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.67%  '
$ws.Range('E9').Value = '  +2.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0624'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.19%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').Value = '1.922.46'
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').Value = '1.675.97'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('E14').Value = '  +2.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.558'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.84'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').Value = '27.177.85'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '236.02'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.93'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.06%  '
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('E22').Value = '  +1.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.55'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.96%  '
$ws.Range('E24').Value = '  -1.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.89'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.41'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  +1.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.17'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('E32').Value = '  +0.54%  '
$ws.Range('D33').Value = '1.549.21'
$ws.Range('E33').Value = '  +1.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.605'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.70%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.948'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.39'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.46%  '
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.07'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.54%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '69.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.84%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = '1.828.69'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '89.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('E48').Value = '  +4.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.22'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.12%  '
$ws.Range('E51').Value = '  +0.12%  '
